$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string values referenced by the new rows.
$fillIn = "Fill in the required information completely"
$convert = "Convert amount to dollars"

# Rows 5-7 use the "Fill in the required information completely" label,
# rows 8-10 use "Convert amount to dollars"; columns B and C repeat the
# existing "PASSED" / "chrome" values used throughout the sheet.
$ws.Range("A5").Value = $fillIn
$ws.Range("B5").Value = "PASSED"
$ws.Range("C5").Value = "chrome"

$ws.Range("A6").Value = $fillIn
$ws.Range("B6").Value = "PASSED"
$ws.Range("C6").Value = "chrome"

$ws.Range("A7").Value = $fillIn
$ws.Range("B7").Value = "PASSED"
$ws.Range("C7").Value = "chrome"

$ws.Range("A8").Value = $convert
$ws.Range("B8").Value = "PASSED"
$ws.Range("C8").Value = "chrome"

$ws.Range("A9").Value = $convert
$ws.Range("B9").Value = "PASSED"
$ws.Range("C9").Value = "chrome"

$ws.Range("A10").Value = $convert
$ws.Range("B10").Value = "PASSED"
$ws.Range("C10").Value = "chrome"
